# Apply the recorded edits to simulation/scenarios.xlsx
#
# Semantic changes described by the diff:
#   1. On the "sampling" sheet, row 2 (Scenario 0):
#        - G2 (the 2028 effort column) changes from 1 to 0
#        - H2 (Description) changes from "100% effort 2023-2028"
#          to "100% effort 2023-2027"
#          (I2 = SUM(B2:G2) recalculates automatically from 6 to 5)
#   2. The active/selected sheet changes from "lethality" to "sampling",
#      and the selection on the sampling sheet moves from H9 to H11.

$wb = $excel.ActiveWorkbook

$sampling = $wb.Worksheets.Item("sampling")

# --- Update the scenario 0 row ---
$sampling.Range("G2").Value = 0
$sampling.Range("H2").Value = "100% effort 2023-2027"

# --- Update which sheet/cell is active & selected ---
# Activating "sampling" makes it the workbook's active tab and clears
# tabSelected from whichever sheet previously had it (lethality).
$sampling.Activate()
$sampling.Range("H11").Select()
